# Apply updated "People" counts for Country "C1" rows (per diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 40
$ws.Range("E6").Value = 10
$ws.Range("E7").Value = 14
$ws.Range("E14").Value = 10
$ws.Range("E15").Value = 4
$ws.Range("E16").Value = 12
$ws.Range("E17").Value = 4
$ws.Range("E18").Value = 12
$ws.Range("E19").Value = 18
$ws.Range("E26").Value = 6
$ws.Range("E27").Value = 16
$ws.Range("E28").Value = 2
$ws.Range("E29").Value = 12
$ws.Range("E30").Value = 4
$ws.Range("E31").Value = 12
$ws.Range("E38").Value = 20
$ws.Range("E39").Value = 8
$ws.Range("E40").Value = 10
$ws.Range("E41").Value = 8
$ws.Range("E42").Value = 6
$ws.Range("E43").Value = 18
$ws.Range("E50").Value = 6
$ws.Range("E51").Value = 10
$ws.Range("E52").Value = 20
$ws.Range("E53").Value = 20
$ws.Range("E54").Value = 2
$ws.Range("E55").Value = 4
